# Hide 3336 in service portfolio
#
# The "Express Same Day" service (service code 3336, request code
# EXPRESS_NORDIC_SAME_DAY) is removed from the "Booking & SG API" sheet.
# It currently lives in row 52; deleting the whole row shifts every
# subsequent row up by one and shrinks the sheet's used range / filter
# range by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire row for the "Express Same Day" / 3336 / EXPRESS_NORDIC_SAME_DAY service.
$ws.Rows(52).Delete()

# The AutoFilter range does not shrink automatically when rows are deleted,
# so reapply it over the new (one-row-shorter) data range.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:P61").AutoFilter()

# The workbook-level hidden "_FilterDatabase" defined name also needs to be
# brought back in sync with the new AutoFilter range.
$filterDatabaseName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterDatabaseName.RefersTo = "='Booking & SG API'!`$A`$1:`$P`$61"
